# 9.5.2 workbook update:
#  - B1 header text (Russian) is replaced with a new translation.
#  - A new "2023" data column (Q) is appended, copying the formatting of
#    the existing "2022" column (P) for both the year-header row and the
#    data row.
#  - Selection is reset back to the top-left cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the Russian header text in B1 -------------------------------
$ws.Range("B1").Value = "9.5.2 Количество исследователей (в эквиваленте полной занятости) на миллион жителей"

# --- 2. Add the new 2023 column (Q), mirroring column P's formatting -------
$pYear = $ws.Cells.Item(4, 16)   # P4 = 2022
$qYear = $ws.Cells.Item(4, 17)   # Q4 = new 2023
$pYear.Copy() | Out-Null
$qYear.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$qYear.Value = 2023

$pData = $ws.Cells.Item(5, 16)   # P5 = 605
$qData = $ws.Cells.Item(5, 17)   # Q5 = new 631
$pData.Copy() | Out-Null
$qData.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$qData.Value = 631

$excel.CutCopyMode = 0

# --- 3. Reset selection back to A1 ------------------------------------------
$ws.Range("A1").Select() | Out-Null
